$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-151 down to 51-152.
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new weekly data point.
$ws.Range("A50").Value = 4
$ws.Range("B50").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C50").Value = "Los Lagos"
$ws.Range("D50").Value = 44544
$ws.Range("E50").Value = 10
$ws.Range("F50").Value = 100112039
$ws.Range("G50").Value = "Ciboulette"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 280
$ws.Range("K50").Value = 2500
$ws.Range("L50").Value = 2500
$ws.Range("M50").Value = 2500
$ws.Range("N50").Value = '$/docena de atados'
$ws.Range("O50").Value = "Región Metropolitana"
$ws.Range("P50").Value = 833
$ws.Range("Q50").Value = 3
$ws.Range("R50").Value = "Hortaliza"
